{"js": "const replacements = [\n  [\"Ativa\u00e7\u00e3o: 01/01/2018\", \"Ativa\u00e7\u00e3o: 01/01/2024\"],\n  [\"Apresentar os princ\u00edpios b\u00e1sicos de Ci\u00eancias dos Materiais, destacando a correla\u00e7\u00e3o entre o comportamento mec\u00e2nico dos metais e os aspectos microestruturais, para aplica\u00e7\u00e3o em Engenharia.\", \"Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais e t\u00eam como objetivo gerar compet\u00eancias no desenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de materiais e a redu\u00e7\u00e3o de ocorr\u00eancia de falhas estruturais. Para tanto, a disciplina estabelece correla\u00e7\u00f5es com outras do curso de Engenharia de Materiais como LOM3013 \u2013 Ci\u00eancia dos Materiais, LOM3057 \u2013 Introdu\u00e7\u00e3o aos Materiais Polim\u00e9ricos, LOM3032 - Cer\u00e2mica F\u00edsica e LOM3011- Ensaios Mec\u00e2nicos. Desta forma, s\u00e3o apresentadas a correla\u00e7\u00e3o entre propriedades e microestrutura de materiais para aplica\u00e7\u00f5es em Engenharia permitindo aos alunos a pr\u00e1tica da reda\u00e7\u00e3o cient\u00edfica e da busca bibliogr\u00e1fica para incentivar a solu\u00e7\u00e3o de problemas em engenharia.\"],\n  [\"1. Introdu\u00e7\u00e3o ao conceito de propriedades mec\u00e2nicas.2. Deforma\u00e7\u00e3o pl\u00e1stica de monocristais e policristais.3. Teoria das discord\u00e2ncias.4. Mecanismos de endurecimento. 5. Comportamento mec\u00e2nico dos materiais met\u00e1licos6. Influ\u00eancias ambientais e t\u00e9rmicas no comportamento mec\u00e2nico. An\u00e1lise de falhas.\", \"1. Introdu\u00e7\u00e3o ao conceito de propriedades mec\u00e2nicas. 2. Elasticidade e Mecanismos de deforma\u00e7\u00e3o pl\u00e1stica. 3. Teoria das discord\u00e2ncias. 4.Mecanismos de endurecimento. 5. Comportamento mec\u00e2nico dos materiais met\u00e1licos. 6. Estudo comparativo de propriedades mec\u00e2nicas de materiais met\u00e1licos, cer\u00e2micos e polim\u00e9ricos. 7. Influ\u00eancia da temperatura no comportamento mec\u00e2nico de materiais. 8. Introdu\u00e7\u00e3o b\u00e1sica \u00e0 an\u00e1lise de falhas de materiais d\u00facteis e fr\u00e1geis.\"],\n  [\"1.INTRODU\u00c7\u00c3O AO CONCEITO DE PROPRIEDADES MEC\u00c2NICAS: Conceitos e rela\u00e7\u00f5es entre microestrutura e propriedades mec\u00e2nicas. Comportamento el\u00e1stico e pl\u00e1stico de metais e ligas. Rela\u00e7\u00f5es entre tens\u00e3o e deforma\u00e7\u00e3o uniaxiais para regime pl\u00e1stico.2.DEFORMA\u00c7\u00c3O PL\u00c1STICA DE MONOCRISTAIS E POLICRISTAIS: Deforma\u00e7\u00e3o pl\u00e1stica e encruamento de monocristais. Sistemas de deslizamento. Deforma\u00e7\u00e3o por macla\u00e7\u00e3o e movimenta\u00e7\u00e3o de discord\u00e2ncias. Movimento relativo de gr\u00e3os.3.TEORIA DAS DISCORD\u00c2NCIAS: Classifica\u00e7\u00e3o, observa\u00e7\u00e3o e fontes de discord\u00e2ncias. Multiplica\u00e7\u00e3o e intera\u00e7\u00e3o de discord\u00e2ncias. For\u00e7as entre discord\u00e2ncias. For\u00e7as atuantes sobre discord\u00e2ncias. Campos de tens\u00e3o e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discord\u00e2ncias. Subestruturas de discord\u00e2ncias. 4.MECANISMOS DE ENDURECIMENTO: Endurecimento por deforma\u00e7\u00e3o pl\u00e1stica: Encruamento. Aumento da resist\u00eancia devido aos contornos de gr\u00e3o e \u00e0 forma\u00e7\u00e3o de c\u00e9lulas e subgr\u00e3os. Rela\u00e7\u00e3o de Hall-Petch. Endurecimento por solu\u00e7\u00e3o s\u00f3lida. Endurecimento por precipita\u00e7\u00e3o. Diagrama Ferro-Carbono. Curvas TTT. A\u00e7os comuns e especiais. Tratamentos t\u00e9rmicos em a\u00e7os; Transforma\u00e7\u00e3o martens\u00edtica.5.COMPORTAMENTO MEC\u00c2NICO DOS MATERIAIS MET\u00c1LICOS: Rela\u00e7\u00e3o entre microestrutura e propriedades. An\u00e1lise das propriedades em fun\u00e7\u00e3o de solicita\u00e7\u00f5es est\u00e1ticas e c\u00edclicas. Propriedades em tra\u00e7\u00e3o uniaxial, flu\u00eancia, fadiga de alto ciclo e propaga\u00e7\u00e3o de trincas por fadiga. Impacto e a transi\u00e7\u00e3o d\u00factil-fr\u00e1gil.6.Influ\u00eancias ambientais e t\u00e9rmicas sobre o comportamento mec\u00e2nico dos metais. An\u00e1lise de falhas em componentes.\", \"1.INTRODU\u00c7\u00c3O AO CONCEITO DE PROPRIEDADES MEC\u00c2NICAS: Conceitos e rela\u00e7\u00f5es entre microestrutura e propriedades mec\u00e2nicas de materiais. Comportamento el\u00e1stico e pl\u00e1stico de metais e ligas. 2. MECANISMOS DE DEFORMA\u00c7\u00c3O PL\u00c1STICA: Sistemas de deslizamento e movimenta\u00e7\u00e3o de discord\u00e2ncias. Deforma\u00e7\u00e3o por macla\u00e7\u00e3o Movimento relativo de gr\u00e3os. Difus\u00e3o. 3. TEORIA DAS DISCORD\u00c2NCIAS: Classifica\u00e7\u00e3o, observa\u00e7\u00e3o e fontes de discord\u00e2ncias. Multiplica\u00e7\u00e3o e intera\u00e7\u00e3o de discord\u00e2ncias. For\u00e7as entre discord\u00e2ncias. For\u00e7as atuantes sobre discord\u00e2ncias. Campos de tens\u00e3o e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discord\u00e2ncias. Subestruturas de discord\u00e2ncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deforma\u00e7\u00e3o pl\u00e1stica: Encruamento. Aumento da resist\u00eancia devido aos contornos de gr\u00e3o. Rela\u00e7\u00e3o de Hall-Petch. Endurecimento por solu\u00e7\u00e3o s\u00f3lida. Endurecimento por precipita\u00e7\u00e3o. A\u00e7os comuns e especiais. Tratamentos t\u00e9rmicos em a\u00e7os. 5. COMPORTAMENTO MEC\u00c2NICO DOS MATERIAIS MET\u00c1LICOS: Rela\u00e7\u00e3o entre microestrutura e propriedades. An\u00e1lise das propriedades em fun\u00e7\u00e3o de solicita\u00e7\u00f5es est\u00e1ticas e c\u00edclicas. Propriedades em tra\u00e7\u00e3o uniaxial, flu\u00eancia, fadiga de alto ciclo e propaga\u00e7\u00e3o de trincas por fadiga. Impacto e a transi\u00e7\u00e3o d\u00factil-fr\u00e1gil. 6. COMPORTAMENTO MEC\u00c2NICO DE MATERIAIS CER\u00c2MICOS E POLIM\u00c9RICOS: Estudo comparativo de propriedades mec\u00e2nicas de materiais met\u00e1licos, cer\u00e2micos e polim\u00e9ricos 7. Influ\u00eancia da temperatura sobre o comportamento mec\u00e2nico de materiais. Aspectos b\u00e1sicos  da  an\u00e1lise de falhas em materiais met\u00e1licos, cer\u00e2micos e polim\u00e9ricos.\"],\n  [\"Este curso dever\u00e1 conter duas avalia\u00e7\u00f5es escritas denominadas P1 e P2. A P2 dever\u00e1 englobar toda a mat\u00e9ria ministrada ao longo do semestre, abrangendo todos os t\u00f3picos previstos na ementa.\", \"Os alunos ser\u00e3o avaliados quanto \u00e0s habilidades gerais em fun\u00e7\u00e3o da participa\u00e7\u00e3o ativa nas aulas. Ser\u00e3o realizadas duas provas escritas P1 e P2, lista de exerc\u00edcios (E) e/ou monografias (M).\"],\n  [\"A m\u00e9dia do semestre ser\u00e1 computada com base na rela\u00e7\u00e3o:M=(P1+2P2)/3\", \"A nota final (NF) do semestre ser\u00e1 calculada pela express\u00e3o: NF = [(P1 + P2)/2] x 0,9 + (E e/ou M) x 0,1. Em caso de aplica\u00e7\u00e3o de Exerc\u00edcios (E), ou prepara\u00e7\u00e3o de monografias (M) e /ou E e M, ser\u00e1 determinada a m\u00e9dia aritm\u00e9tica entre as notas e multiplicadas pelo fator 0,1 para o c\u00e1lculo de NF.\"],\n  [\"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre.A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 computada com base na rela\u00e7\u00e3o abaixo:MF=(M+RC)/2\", \"1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 2. A. S. Lisb\u00e3o, Estrutura e propriedades dos pol\u00edmeros, EduFSCar, S\u00e3o Carlos, 2009. 3. T. H. Courtney, Mechanical Behavior of Materials, Waveland Press, 2005. 4. A. K. Bhargava, Engineering Materials: Polymers, Ceramics and Composites, PHI Learning Pvt. Ltd., 2012. 5.Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 2007. 6. Hull, D. Introduction to Dislocations, Pergamon Press, 1965. 7. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967. 8. Reed-Hill, R.E. Princ\u00edpios de Metalurgia F\u00edsica, Ed. Guanabara Dois, 1982. 9. Van Vlack, L.H. Princ\u00edpios de Ci\u00eancia dos Materiais, Ed. Edgard Blucher Ltda., 1970. 10. Costa e Silva, A. L., Mei, P. R. A\u00e7os e Ligas especiais, Ed. Edgar Bl\u00fccher, 2008. 11. Dieter, G.E. Metalurgia Mec\u00e2nica, Ed. Guanabara Dois, 1986.  12. Callister, W. Ci\u00eancia e engenharia dos materiais: Uma introdu\u00e7\u00e3o, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2008. 13. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993.\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match for: \" + oldText.slice(0, 40) + \" \u2014 got \" + results.items.length);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"Ativa\u00e7\u00e3o: 01/01/2018\", \"Ativa\u00e7\u00e3o: 01/01/2024\"),\n  @(\"Apresentar os princ\u00edpios b\u00e1sicos de Ci\u00eancias dos Materiais, destacando a correla\u00e7\u00e3o entre o comportamento mec\u00e2nico dos metais e os aspectos microestruturais, para aplica\u00e7\u00e3o em Engenharia.\", \"Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais e t\u00eam como objetivo gerar compet\u00eancias no desenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de materiais e a redu\u00e7\u00e3o de ocorr\u00eancia de falhas estruturais. Para tanto, a disciplina estabelece correla\u00e7\u00f5es com outras do curso de Engenharia de Materiais como LOM3013 \u2013 Ci\u00eancia dos Materiais, LOM3057 \u2013 Introdu\u00e7\u00e3o aos Materiais Polim\u00e9ricos, LOM3032 - Cer\u00e2mica F\u00edsica e LOM3011- Ensaios Mec\u00e2nicos. Desta forma, s\u00e3o apresentadas a correla\u00e7\u00e3o entre propriedades e microestrutura de materiais para aplica\u00e7\u00f5es em Engenharia permitindo aos alunos a pr\u00e1tica da reda\u00e7\u00e3o cient\u00edfica e da busca bibliogr\u00e1fica para incentivar a solu\u00e7\u00e3o de problemas em engenharia.\"),\n  @(\"1. Introdu\u00e7\u00e3o ao conceito de propriedades mec\u00e2nicas.2. Deforma\u00e7\u00e3o pl\u00e1stica de monocristais e policristais.3. Teoria das discord\u00e2ncias.4. Mecanismos de endurecimento. 5. Comportamento mec\u00e2nico dos materiais met\u00e1licos6. Influ\u00eancias ambientais e t\u00e9rmicas no comportamento mec\u00e2nico. An\u00e1lise de falhas.\", \"1. Introdu\u00e7\u00e3o ao conceito de propriedades mec\u00e2nicas. 2. Elasticidade e Mecanismos de deforma\u00e7\u00e3o pl\u00e1stica. 3. Teoria das discord\u00e2ncias. 4.Mecanismos de endurecimento. 5. Comportamento mec\u00e2nico dos materiais met\u00e1licos. 6. Estudo comparativo de propriedades mec\u00e2nicas de materiais met\u00e1licos, cer\u00e2micos e polim\u00e9ricos. 7. Influ\u00eancia da temperatura no comportamento mec\u00e2nico de materiais. 8. Introdu\u00e7\u00e3o b\u00e1sica \u00e0 an\u00e1lise de falhas de materiais d\u00facteis e fr\u00e1geis.\"),\n  @(\"1.INTRODU\u00c7\u00c3O AO CONCEITO DE PROPRIEDADES MEC\u00c2NICAS: Conceitos e rela\u00e7\u00f5es entre microestrutura e propriedades mec\u00e2nicas. Comportamento el\u00e1stico e pl\u00e1stico de metais e ligas. Rela\u00e7\u00f5es entre tens\u00e3o e deforma\u00e7\u00e3o uniaxiais para regime pl\u00e1stico.2.DEFORMA\u00c7\u00c3O PL\u00c1STICA DE MONOCRISTAIS E POLICRISTAIS: Deforma\u00e7\u00e3o pl\u00e1stica e encruamento de monocristais. Sistemas de deslizamento. Deforma\u00e7\u00e3o por macla\u00e7\u00e3o e movimenta\u00e7\u00e3o de discord\u00e2ncias. Movimento relativo de gr\u00e3os.3.TEORIA DAS DISCORD\u00c2NCIAS: Classifica\u00e7\u00e3o, observa\u00e7\u00e3o e fontes de discord\u00e2ncias. Multiplica\u00e7\u00e3o e intera\u00e7\u00e3o de discord\u00e2ncias. For\u00e7as entre discord\u00e2ncias. For\u00e7as atuantes sobre discord\u00e2ncias. Campos de tens\u00e3o e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discord\u00e2ncias. Subestruturas de discord\u00e2ncias. 4.MECANISMOS DE ENDURECIMENTO: Endurecimento por deforma\u00e7\u00e3o pl\u00e1stica: Encruamento. Aumento da resist\u00eancia devido aos contornos de gr\u00e3o e \u00e0 forma\u00e7\u00e3o de c\u00e9lulas e subgr\u00e3os. Rela\u00e7\u00e3o de Hall-Petch. Endurecimento por solu\u00e7\u00e3o s\u00f3lida. Endurecimento por precipita\u00e7\u00e3o. Diagrama Ferro-Carbono. Curvas TTT. A\u00e7os comuns e especiais. Tratamentos t\u00e9rmicos em a\u00e7os; Transforma\u00e7\u00e3o martens\u00edtica.5.COMPORTAMENTO MEC\u00c2NICO DOS MATERIAIS MET\u00c1LICOS: Rela\u00e7\u00e3o entre microestrutura e propriedades. An\u00e1lise das propriedades em fun\u00e7\u00e3o de solicita\u00e7\u00f5es est\u00e1ticas e c\u00edclicas. Propriedades em tra\u00e7\u00e3o uniaxial, flu\u00eancia, fadiga de alto ciclo e propaga\u00e7\u00e3o de trincas por fadiga. Impacto e a transi\u00e7\u00e3o d\u00factil-fr\u00e1gil.6.Influ\u00eancias ambientais e t\u00e9rmicas sobre o comportamento mec\u00e2nico dos metais. An\u00e1lise de falhas em componentes.\", \"1.INTRODU\u00c7\u00c3O AO CONCEITO DE PROPRIEDADES MEC\u00c2NICAS: Conceitos e rela\u00e7\u00f5es entre microestrutura e propriedades mec\u00e2nicas de materiais. Comportamento el\u00e1stico e pl\u00e1stico de metais e ligas. 2. MECANISMOS DE DEFORMA\u00c7\u00c3O PL\u00c1STICA: Sistemas de deslizamento e movimenta\u00e7\u00e3o de discord\u00e2ncias. Deforma\u00e7\u00e3o por macla\u00e7\u00e3o Movimento relativo de gr\u00e3os. Difus\u00e3o. 3. TEORIA DAS DISCORD\u00c2NCIAS: Classifica\u00e7\u00e3o, observa\u00e7\u00e3o e fontes de discord\u00e2ncias. Multiplica\u00e7\u00e3o e intera\u00e7\u00e3o de discord\u00e2ncias. For\u00e7as entre discord\u00e2ncias. For\u00e7as atuantes sobre discord\u00e2ncias. Campos de tens\u00e3o e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discord\u00e2ncias. Subestruturas de discord\u00e2ncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deforma\u00e7\u00e3o pl\u00e1stica: Encruamento. Aumento da resist\u00eancia devido aos contornos de gr\u00e3o. Rela\u00e7\u00e3o de Hall-Petch. Endurecimento por solu\u00e7\u00e3o s\u00f3lida. Endurecimento por precipita\u00e7\u00e3o. A\u00e7os comuns e especiais. Tratamentos t\u00e9rmicos em a\u00e7os. 5. COMPORTAMENTO MEC\u00c2NICO DOS MATERIAIS MET\u00c1LICOS: Rela\u00e7\u00e3o entre microestrutura e propriedades. An\u00e1lise das propriedades em fun\u00e7\u00e3o de solicita\u00e7\u00f5es est\u00e1ticas e c\u00edclicas. Propriedades em tra\u00e7\u00e3o uniaxial, flu\u00eancia, fadiga de alto ciclo e propaga\u00e7\u00e3o de trincas por fadiga. Impacto e a transi\u00e7\u00e3o d\u00factil-fr\u00e1gil. 6. COMPORTAMENTO MEC\u00c2NICO DE MATERIAIS CER\u00c2MICOS E POLIM\u00c9RICOS: Estudo comparativo de propriedades mec\u00e2nicas de materiais met\u00e1licos, cer\u00e2micos e polim\u00e9ricos 7. Influ\u00eancia da temperatura sobre o comportamento mec\u00e2nico de materiais. Aspectos b\u00e1sicos  da  an\u00e1lise de falhas em materiais met\u00e1licos, cer\u00e2micos e polim\u00e9ricos.\"),\n  @(\"Este curso dever\u00e1 conter duas avalia\u00e7\u00f5es escritas denominadas P1 e P2. A P2 dever\u00e1 englobar toda a mat\u00e9ria ministrada ao longo do semestre, abrangendo todos os t\u00f3picos previstos na ementa.\", \"Os alunos ser\u00e3o avaliados quanto \u00e0s habilidades gerais em fun\u00e7\u00e3o da participa\u00e7\u00e3o ativa nas aulas. Ser\u00e3o realizadas duas provas escritas P1 e P2, lista de exerc\u00edcios (E) e/ou monografias (M).\"),\n  @(\"A m\u00e9dia do semestre ser\u00e1 computada com base na rela\u00e7\u00e3o:M=(P1+2P2)/3\", \"A nota final (NF) do semestre ser\u00e1 calculada pela express\u00e3o: NF = [(P1 + P2)/2] x 0,9 + (E e/ou M) x 0,1. Em caso de aplica\u00e7\u00e3o de Exerc\u00edcios (E), ou prepara\u00e7\u00e3o de monografias (M) e /ou E e M, ser\u00e1 determinada a m\u00e9dia aritm\u00e9tica entre as notas e multiplicadas pelo fator 0,1 para o c\u00e1lculo de NF.\"),\n  @(\"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre.A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 computada com base na rela\u00e7\u00e3o abaixo:MF=(M+RC)/2\", \"1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 2. A. S. Lisb\u00e3o, Estrutura e propriedades dos pol\u00edmeros, EduFSCar, S\u00e3o Carlos, 2009. 3. T. H. Courtney, Mechanical Behavior of Materials, Waveland Press, 2005. 4. A. K. Bhargava, Engineering Materials: Polymers, Ceramics and Composites, PHI Learning Pvt. Ltd., 2012. 5.Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 2007. 6. Hull, D. Introduction to Dislocations, Pergamon Press, 1965. 7. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967. 8. Reed-Hill, R.E. Princ\u00edpios de Metalurgia F\u00edsica, Ed. Guanabara Dois, 1982. 9. Van Vlack, L.H. Princ\u00edpios de Ci\u00eancia dos Materiais, Ed. Edgard Blucher Ltda., 1970. 10. Costa e Silva, A. L., Mei, P. R. A\u00e7os e Ligas especiais, Ed. Edgar Bl\u00fccher, 2008. 11. Dieter, G.E. Metalurgia Mec\u00e2nica, Ed. Guanabara Dois, 1986.  12. Callister, W. Ci\u00eancia e engenharia dos materiais: Uma introdu\u00e7\u00e3o, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2008. 13. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993.\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Forward = $true\n  $find.Wrap = 0\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n  $ok = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 1)\n  if (-not $ok) {\n    throw \"Find/Replace failed for: \" + $pair[0].Substring(0, [Math]::Min(40, $pair[0].Length))\n  }\n}\n"}
